# Apply the ValueSet-attribution-source update:
#  - Version bumped 5.0.0 -> 6.0.0
#  - Date bumped to the new publish timestamp
#  - Publisher value filled in ("Alvearie Team")
#  - The old duplicated "Contact" / "No display for ContactDetail" row is
#    replaced by a single "Jurisdiction" / "United States of America" row
#    (net: one row removed from the sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B3").Value = "6.0.0"
$ws.Range("B8").Value = "2022-01-21T20:46:54+00:00"
$ws.Range("B9").Value = "Alvearie Team"
$ws.Range("A10").Value = "Jurisdiction"
$ws.Range("B10").Value = "United States of America"

# Remove the now-redundant second "Contact" row (old row 11); everything
# below shifts up by one, collapsing the sheet from 15 rows to 14.
$ws.Rows(11).Delete()
